$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '29.322.76'
$ws.Range("E2").Value = '  -0.30%  '

$ws.Range("D3").Value = '1.845.23'
$ws.Range("E3").Value = '  -0.24%  '

Set-TextValue $ws.Range("D4") '0.9976'
$ws.Range("E4").Value = '  -0.22%  '

Set-TextValue $ws.Range("D5") '239.88'
$ws.Range("E5").Value = '  -0.36%  '

Set-TextValue $ws.Range("D6") '0.6262'
$ws.Range("E6").Value = '  -0.64%  '

Set-TextValue $ws.Range("D7") '0.9974'
$ws.Range("E7").Value = '  -0.32%  '

Set-TextValue $ws.Range("D8") '0.07604'
$ws.Range("E8").Value = '  -1.30%  '

Set-TextValue $ws.Range("D9") '0.2898'
$ws.Range("E9").Value = '  -1.46%  '

Set-TextValue $ws.Range("D10") '24.72'
$ws.Range("E10").Value = '  +0.76%  '

Set-TextValue $ws.Range("D11") '0.07725'
$ws.Range("E11").Value = '  -0.32%  '

Set-TextValue $ws.Range("D12") '5.018'
$ws.Range("E12").Value = '  -0.22%  '

Set-TextValue $ws.Range("D13") '0.6777'
$ws.Range("E13").Value = '  -0.39%  '

Set-TextValue $ws.Range("D14") '0.00001048'
$ws.Range("E14").Value = '  -3.71%  '

Set-TextValue $ws.Range("D15") '82.92'
$ws.Range("E15").Value = '  -0.90%  '

$ws.Range("E16").Value = '  -0.28%  '

$ws.Range("D17").Value = '29.352.19'
$ws.Range("E17").Value = '  -0.34%  '

Set-TextValue $ws.Range("D18") '227.66'
$ws.Range("E18").Value = '  -0.62%  '

Set-TextValue $ws.Range("D19") '12.31'
$ws.Range("E19").Value = '  -1.28%  '

Set-TextValue $ws.Range("D20") '0.9972'
$ws.Range("E20").Value = '  -0.32%  '

Set-TextValue $ws.Range("D21") '7.452'
$ws.Range("E21").Value = '  +0.03%  '

Set-TextValue $ws.Range("D22") '0.9985'
$ws.Range("E22").Value = '  -0.23%  '

Set-TextValue $ws.Range("D23") '158.33'
$ws.Range("E23").Value = '  +0.68%  '

Set-TextValue $ws.Range("D24") '0.1382'
$ws.Range("E24").Value = '  -0.50%  '

Set-TextValue $ws.Range("D25") '8.408'
$ws.Range("E25").Value = '  +0.62%  '

$ws.Range("E26").Value = '  -0.25%  '

Set-TextValue $ws.Range("D27") '1.400'
$ws.Range("E27").Value = '  +6.49%  '

$ws.Range("E28").Value = '  -0.68%  '

Set-TextValue $ws.Range("D29") '0.05595'
$ws.Range("E29").Value = '  -1.46%  '

Set-TextValue $ws.Range("D30") '4.102'
$ws.Range("E30").Value = '  -0.20%  '

Set-TextValue $ws.Range("D31") '4.054'
$ws.Range("E31").Value = '  +0.11%  '

Set-TextValue $ws.Range("D32") '1.160'
$ws.Range("E32").Value = '  +0.12%  '

Set-TextValue $ws.Range("D33") '1.827'
$ws.Range("E33").Value = '  -1.27%  '

Set-TextValue $ws.Range("D34") '0.6949'
$ws.Range("E34").Value = '  -1.95%  '

Set-TextValue $ws.Range("D35") '2.581'
$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").Value = '1.223.52'
$ws.Range("E37").Value = '  -0.54%  '

Set-TextValue $ws.Range("D38") '2.717'
$ws.Range("E38").Value = '  -2.22%  '

Set-TextValue $ws.Range("D39") '6.342'
$ws.Range("E39").Value = '  -2.40%  '

Set-TextValue $ws.Range("D40") '0.9007'
$ws.Range("E40").Value = '  -1.55%  '

Set-TextValue $ws.Range("D41") '0.9970'
$ws.Range("E41").Value = '  -0.35%  '

Set-TextValue $ws.Range("D42") '101.19'
$ws.Range("E42").Value = '  -0.29%  '

Set-TextValue $ws.Range("D43") '65.44'
$ws.Range("E43").Value = '  -1.30%  '

Set-TextValue $ws.Range("D44") '7.181'
$ws.Range("E44").Value = '  +0.24%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D45") '0.00000000118'
$ws.Range("E45").Value = '  -2.79%  '

$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D46") '0.3986'
$ws.Range("E46").Value = '  -0.75%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D47") '8.974'
$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D48") '1.677'
$ws.Range("E48").Value = '  -0.72%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D49") '0.1137'
$ws.Range("E49").Value = '  +1.15%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.05698'
$ws.Range("E50").Value = '  -0.25%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D51") '0.4617'
$ws.Range("E51").Value = '  -0.23%  '
